$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Gemma-7B-Instruct"
$ws.Range("B8").Value = "0.81 ± 0.39"
$ws.Range("C8").Value = "0.42 ± 0.79"
$ws.Range("D8").Value = "0.26 ± 0.56"
$ws.Range("E8").Value = "0.01 ± 0.02"
$ws.Range("F8").Value = "0.14 ± 0.09"
$ws.Range("G8").Value = "0.03 ± 0.04"
$ws.Range("H8").Value = "0.12 ± 0.08"
$ws.Range("I8").Value = "0.15 ± 0.1"
$ws.Range("J8").Value = "0.71 ± 0.3"
$ws.Range("K8").Value = "0.73 ± 0.31"
$ws.Range("L8").Value = "0.72 ± 0.3"
$ws.Range("M8").Value = "0.68 ± 0.29"
$ws.Range("N8").Value = "0.83 ± 0.35"
$ws.Range("O8").Value = "0.1 ± 0.09"
$ws.Range("P8").Value = "0.53 ± 0.25"
$ws.Range("Q8").Value = "9.68 ± 1.44"
$ws.Range("R8").Value = "0.15 ± 0.00"
$ws.Range("S8").Value = "0.77 ± 0.33"
$ws.Range("T8").Value = "0.85 ± 0.36"
$ws.Range("U8").Value = "3.63 ± 1.94"
$ws.Range("V8").Value = "0.81 ± 0.38"
$ws.Range("W8").Value = "0.77 ± 0.33"
$ws.Range("X8").Value = "1.37 ± 0.62"
